# Edit described by commit "Changes in naming of csv files":
#   - slide 3 ("Data preparation") body placeholder:
#       * "all_pressreleases.csv" -> "pressreleases_all.csv"
#       * "all_speaches.csv"      -> "speeches_all.csv"
#     (PowerPoint splits the edited paragraphs into several runs around
#     the retyped words; we reproduce the same run boundaries.)
#   - notes page of slide 3: nudge the slide-image placeholder's left
#     offset by a hair (381300 EMU -> 381000 EMU).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)

# --- Paragraph 3: "Number of press releases: 3624 (all_pressreleases.csv)"
#     -> "Number of press releases: 3624 (pressreleases_all.csv)"
# Split off "(all_pressreleases.csv)" into its own run, retyped.
$para3 = $sh.TextFrame.TextRange.Paragraphs(3)
$sel3a = $para3.Characters(32, 24)
$sel3a.Text = "(pressreleases_all.csv)"

# Carve the trailing ")" into its own run.
$para3b = $sh.TextFrame.TextRange.Paragraphs(3)
$sel3b = $para3b.Characters($para3b.Length - 1, 1)
$sel3b.Text = ")"

# --- Paragraph 4: "Number of speeches: 809 (all_speaches.csv)"
#     -> "Number of speeches: 809 (speeches_all.csv)"
# Split off "(all_speaches.csv)" into its own run, retyped.
$para4 = $sh.TextFrame.TextRange.Paragraphs(4)
$sel4a = $para4.Characters(25, 18)
$sel4a.Text = "(speeches_all.csv)"

# Carve the trailing ")" into its own run.
$para4b = $sh.TextFrame.TextRange.Paragraphs(4)
$sel4b = $para4b.Characters($para4b.Length - 1, 1)
$sel4b.Text = ")"

# Carve "Number of speeches: " / "809 " into separate runs.
$para4c = $sh.TextFrame.TextRange.Paragraphs(4)
$sel4c = $para4c.Characters(21, 4)
$sel4c.Text = "809 "

# --- Notes page: nudge the slide-image placeholder position slightly.
$notesPage = $s.NotesPage
$imgPh = $notesPage.Shapes.Item(1)
$imgPh.Left = 30.0
$imgPh.Top = 54.0
